$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that vary per-row and must be permuted: D,L,M,N,O,P,Q,S,T
# Column index: D=4, L=12, M=13, N=14, O=15, P=16, Q=17, S=19, T=20
$cols = @(4,12,13,14,15,16,17,19,20)

# Snapshot original values for rows 2..20 before overwriting anything
$orig = @{}
for ($r = 2; $r -le 20; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $orig[$r] = $rowVals
}

# Mapping: destination row -> source row (values copied from source row in the snapshot)
$mapping = @{
    2 = 8
    3 = 13
    4 = 2
    5 = 4
    6 = 5
    7 = 7
    8 = 20
    9 = 18
    10 = 6
    11 = 3
    12 = 12
    13 = 16
    14 = 10
    15 = 19
    16 = 15
    17 = 17
    18 = 14
    19 = 9
    20 = 11
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $orig[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value = $srcVals[$c]
    }
}
